# Generate Report for Handback
# Marks the zh-cn and de-de localization rows as handed back: updates the
# "Status" text, fills in the "Latest Target File" (as a hyperlink to the
# source .md) and "Latest Handback File" columns, stamps the handback
# datetime, and widens the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$mdFileName = "dea3305b-635d-45fa-af3a-19e14bc3d44d.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/71496876af611878b6080a5776fb1c56bcb22d66/e2e/dea3305b-635d-45fa-af3a-19e14bc3d44d.md"
$statusText = "Handed back: in sync with en-US"

$zhXlf = "dea3305b-635d-45fa-af3a-19e14bc3d44d.fdb8bbdb6ef1fc41ac4ed20b851a005a639d2c69.zh-cn.xlf"
$deXlf = "dea3305b-635d-45fa-af3a-19e14bc3d44d.fdb8bbdb6ef1fc41ac4ed20b851a005a639d2c69.de-de.xlf"

$zhHandbackTime = "2016-08-24 22:59:27"
$deHandbackTime = "2016-08-24 22:59:34"

# ---- Overview sheet: refresh the per-language status cells ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.16
$overview.Columns.Item(6).ColumnWidth = 29.16

# ---- zh-cn sheet: status + handback file/date ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText
$zh.Hyperlinks.Add($zh.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)
$zh.Range("J2").Value = $zhXlf
$zh.Range("K2").Value = $zhHandbackTime
$zh.Columns.Item(3).ColumnWidth = 29.16
$zh.Columns.Item(9).ColumnWidth = 39.17
$zh.Columns.Item(10).ColumnWidth = 39.17

# ---- de-de sheet: status + handback file/date ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText
$de.Hyperlinks.Add($de.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)
$de.Range("J2").Value = $deXlf
$de.Range("K2").Value = $deHandbackTime
$de.Columns.Item(3).ColumnWidth = 29.16
$de.Columns.Item(9).ColumnWidth = 39.17
$de.Columns.Item(10).ColumnWidth = 39.17
